# ------------------------------------------------------------------
# AAI500_FinalProject_VideoScript.docx edit script
# 1) Replace the single-line "This is our script." paragraph with a
#    multi-paragraph set of recording guidelines.
# 2) Add <w:lastRenderedPageBreak/> to the 2nd and 5th "Israel" table
#    name-cells, and remove it from the 3rd "Mani" name-cell (the page
#    break now falls in a different spot after the text was added above).
# 3) Clean up a few spell/grammar-check proofErr wrappers that Word
#    collapses (and their adjoining runs merge) when the document is
#    re-saved: "No video conferences, no video games.", "...the plant
#    they burn fuel...", and "...the general public...".
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Intro paragraph -> 5 guideline paragraphs -------------------
$introPara = $d.Paragraphs(2)
$introRange = $introPara.Range
if ($introRange.Text -ne "This is our script.`r") {
    Write-Host "WARNING: paragraph 2 did not contain the expected placeholder text"
}
$part1Xml = '<w:p><w:r><w:t>Hey Team! Here are some guides for our video to aim at the professional look and feel I think we want to deliver – I’m open to your suggestions as well.</w:t></w:r></w:p><w:p><w:r><w:t>Try recording this with your cell phones in a horizontal orientation, in selfie mode at about 3 to 4 feet away from you</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Make sure your face and upper part of your body shows in the frame. Try to choose a quiet location and I suggest you sit down to prevent too much uncontrolled motion in the video recording</w:t></w:r><w:r><w:t xml:space="preserve"> (although I was wondering if we should try walking with our phones in selfie mode when applicable)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">I envision the final video as a short documentary to present our research and results, so while </w:t></w:r><w:r><w:t>we</w:t></w:r><w:r><w:t xml:space="preserve"> can smile, picture ourselves presenting to an audience of investors that we’re trying to convince to support our research because it is both profitable and impressive.</w:t></w:r><w:r><w:t xml:space="preserve"> Think of the attitude of someone presenting in a TED talk as well.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">I think for our dress code we should aim for a </w:t></w:r><w:r><w:t xml:space="preserve">dress </w:t></w:r><w:r><w:t>shirt, ideally no hats</w:t></w:r><w:r><w:t xml:space="preserve"> or t-shirts</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> I plan to use only the audio in some parts as a voice-over</w:t></w:r><w:r><w:t>, but I think we should visually include ourselves talking in the video too</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Let me know if you have any questions or suggestions.</w:t></w:r></w:p>'
$introRange.InsertXML($part1Xml)

# --- 2) lastRenderedPageBreak bookkeeping on "Israel" / "Mani" cells -
$israelBreakXml = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Israel</w:t></w:r></w:p>'
$rng = $d.Content
$israelCount = 0
while ($rng.Find.Execute("Israel")) {
    $israelCount = $israelCount + 1
    if ($israelCount -eq 2 -or $israelCount -eq 5) {
        $cellPara = $rng.Paragraphs(1)
        $cellPara.Range.InsertXML($israelBreakXml)
    }
    $rng.Collapse(0)
}
Write-Host "Israel occurrences processed:" $israelCount

$maniPlainXml = '<w:p><w:r><w:t>Mani</w:t></w:r></w:p>'
$rng = $d.Content
$maniCount = 0
while ($rng.Find.Execute("Mani")) {
    $maniCount = $maniCount + 1
    if ($maniCount -eq 3) {
        $cellPara = $rng.Paragraphs(1)
        $cellPara.Range.InsertXML($maniPlainXml)
    }
    $rng.Collapse(0)
}
Write-Host "Mani occurrences processed:" $maniCount

# --- 3) proofErr cleanup paragraphs ----------------------------------
$part3Xml = '<w:p><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t xml:space="preserve">Without it, you couldn’t even watch this video. </w:t></w:r><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>No video conferen</w:t></w:r><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>ces</w:t></w:r><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>, no vi</w:t></w:r><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>deo games.</w:t></w:r><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Electric power is everything. I cannot think of an industry that could live without it today.</w:t></w:r></w:p>'
$rng = $d.Content
$found = $rng.Find.Execute("Without it, you couldn")
if ($found) { $rng.Paragraphs(1).Range.InsertXML($part3Xml) }
else { Write-Host "WARNING: could not locate 'No video conferences' paragraph" }

$part4Xml = '<w:p><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>I</w:t></w:r><w:r><w:t>n the first Stage of the plant they burn fuel and make a turbine spin to induce electricity. That produces a lot of heat.</w:t></w:r></w:p>'
$rng = $d.Content
$found = $rng.Find.Execute("In the first Stage of the")
if ($found) { $rng.Paragraphs(1).Range.InsertXML($part4Xml) }
else { Write-Host "WARNING: could not locate 'first Stage of the plant' paragraph" }

$part7Xml = '<w:p><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t xml:space="preserve">Imagine the impact of such forecast to the demand; the </w:t></w:r><w:r><w:t xml:space="preserve">general </w:t></w:r><w:r><w:t xml:space="preserve">public </w:t></w:r><w:r><w:t xml:space="preserve">as well as companies </w:t></w:r><w:r><w:t>could prepare in advance for their energy needs and make decisions to effectively plan their upcoming energy bill.</w:t></w:r><w:r><w:t xml:space="preserve"> I’m Israel Romero</w:t></w:r></w:p>'
$rng = $d.Content
$found = $rng.Find.Execute("Imagine the impact of such forecast")
if ($found) { $rng.Paragraphs(1).Range.InsertXML($part7Xml) }
else { Write-Host "WARNING: could not locate 'Imagine the impact' paragraph" }

Write-Host "Done."
